# VendAddress.xlsx maintenance edit:
#  - Row 8 (Toronto, AddrID 7) previously had no state abbreviation; set it to "NA".
#  - The Zip4 column (I) was blank for every data row; set it to 0 for rows 2-11
#    so the generated INSERT statements carry an explicit Zip4 value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# St (State) for the Toronto row was blank -> now "NA"
$ws.Range("G8").Value = "NA"

# Zip4 column (I) was blank for all data rows -> now 0
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 9).Value = 0
}
